$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189. Excel's Insert() shifts row 189
# (and everything below it, through the old last row 215) down by one,
# growing the used range from A1:T215 to A1:T216.
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new record.
$ws.Range("A189").Value = 1
$ws.Range("B189").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C189").Value = "Arica y Parinacota"
$ws.Range("D189").Value = 44984
$ws.Range("E189").Value = 15
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100108
$ws.Range("H189").Value = "Tropicales y subtropicales"
$ws.Range("I189").Value = 100108002
$ws.Range("J189").Value = "Mango"
$ws.Range("K189").Value = "Sin especificar"
$ws.Range("L189").Value = "Especial"
$ws.Range("M189").Value = 456
$ws.Range("N189").Value = 4500
$ws.Range("O189").Value = 5000
$ws.Range("P189").Value = 4750
$ws.Range("Q189").Value = "$/bandeja 4 kilos"
$ws.Range("R189").Value = "Perú"
$ws.Range("S189").Value = 1188
$ws.Range("T189").Value = 4
